# Updates league base data, commit: "Atualização de bases das ligas, do dia: 29-05-2024 às 22:54"
# Rows 195-197 (columns B:AD) are cyclically rotated: the old content of row 196
# moves into row 195, the old content of row 197 moves into row 196, and the old
# content of row 195 moves into row 197. Column A (row index) is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hungary NB I")

$data = @(
    @(8209690, "Hungary NB I", 45430.53125, "Paksi", "Kisvarda FC", 2, 1, 1, 0, "H", 1.444, 4.333, 6, 1.45, 4.75, 5, -1.25, 2.025, 1.825, 3.25, 2, 1.85, 0.45, -1, -1, -0.5, 0.4125, -0.5, 0.425),
    @(8209692, "Hungary NB I", 45430.53125, "Puskas Academy", "Debreceni VSC", 4, 1, 2, 1, "H", 1.7, 3.6, 4.5, 1.333, 4.5, 8, -1.5, 2.025, 1.825, 3, 1.875, 1.975, 0.333, -1, -1, 1.025, -1, 0.875, -1),
    @(8209693, "Hungary NB I", 45430.53125, "MOL Fehervar FC", "Diosgyori VTK", 0, 0, 0, 0, "D", 1.571, 4, 5, 1.6, 3.9, 4.75, -1, 2.025, 1.825, 3, 1.925, 1.925, -1, 2.9, -1, -1, 0.825, -1, 0.925)
)

$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD")

$startRow = 195
for ($i = 0; $i -lt $data.Length; $i++) {
    $rowValues = $data[$i]
    $rowNum = $startRow + $i
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $ws.Range($cols[$j] + "$rowNum").Value = $rowValues[$j]
    }
}
